# agendas.xlsx - "mejorando carga de excel y validaciones de dias habiles"
#
# - Corrects the date in B3 (was 27/03/2023, should be 25/03/2023)
# - Appends two new "paddel" schedule rows to the Tabla1 table
#   (dia=24/03/2023 horario=4, dia=26/03/2023 horario=5)
# - Leaves the selection on B6, matching where the user finished typing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- fix the wrong date that had been entered on row 3 ------------------
$ws.Range("B3").Value = 45010

# --- grow the table with the two missing rows ----------------------------
# New rows inherit the look of the row right above them (row 4): same
# fill/border/font for the "tipo_cancha" and "dia" columns, and the
# right-hand border treatment used by the "horario" column elsewhere in
# the table (column C, rows 2/3) for the new "horario" entries.

$row5 = $tbl.ListRows.Add()
$ws.Range("A4:B4").Copy($row5.Range.Cells.Item(1, 1).Resize(1, 2))
$ws.Range("C2").Copy($row5.Range.Cells.Item(1, 3))
$row5.Range.Cells.Item(1, 1).Value = "paddel"
$row5.Range.Cells.Item(1, 2).Value = 45009
$row5.Range.Cells.Item(1, 3).Value = 4

$row6 = $tbl.ListRows.Add()
$ws.Range("A4:B4").Copy($row6.Range.Cells.Item(1, 1).Resize(1, 2))
$ws.Range("C2").Copy($row6.Range.Cells.Item(1, 3))
$row6.Range.Cells.Item(1, 1).Value = "paddel"
$row6.Range.Cells.Item(1, 2).Value = 45011
$row6.Range.Cells.Item(1, 3).Value = 5

# --- match the cursor position left behind after the edit ---------------
$ws.Range("B6").Select() | Out-Null
